$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.872931361198425
$ws.Range("B1").Value = 3.869392871856689
$ws.Range("C1").Value = 2.541174650192261
$ws.Range("D1").Value = 0.9171818494796753
$ws.Range("E1").Value = 0.6013784408569336
